# Generate Report for Handoff
#
# The handoff-status workbook previously tracked a single markdown file
# (ca813ae2-a82f-44ab-bc63-df27ce780301.md). This run regenerates the
# report for a newer commit that hands off an .md file plus the two .png
# images it depends on, so every sheet grows from one data row to three.

$wb = $excel.ActiveWorkbook

# Font color used by the workbook's custom "HyperLink" cell style
# (rgb FF6495ED == RGB(100,149,237) == 15570276). Applied by hand because
# the engine's named-style assignment ($range.Style = "HyperLink") doesn't
# carry the underline/color along with it.
$HyperlinkColor = 15570276

function Style-AsHyperlink($ws, $cellRef) {
    $ws.Range($cellRef).Font.Underline = $true
    $ws.Range($cellRef).Font.Color = $HyperlinkColor
}

function Add-StyledHyperlink($ws, $cellRef, $url, $display) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $display)
    Style-AsHyperlink $ws $cellRef
}

# ---------------------------------------------------------------------
# New source / target file identities for this handoff
# ---------------------------------------------------------------------
$png1Name = "12bf62c6-adf7-4bbf-a2c2-2d4ed70955e8.png"
$png2Name = "8597cde4-2720-4142-aafa-7d6024e663ca.png"
$mdName   = "b74fb850-0273-4207-9667-2a3042d3febb.md"

$png1TargetName = "cd77bedb8fe0827b3a68ba2cf1e18ed7d8c77dba.png"
$png2TargetName = "a208677ae8f632633f73302ae7c80f86309edc76.png"
$mdTargetZhCn   = "b74fb850-0273-4207-9667-2a3042d3febb.fe4c5be2c84d2bdff269b6e14bc9b7aa72791893.zh-cn.xlf"
$mdTargetDeDe   = "b74fb850-0273-4207-9667-2a3042d3febb.fe4c5be2c84d2bdff269b6e14bc9b7aa72791893.de-de.xlf"

$readyStatus  = "Ready for handoff"
$overviewDate = "2016-55-20 16:55:23"

$srcUrlBase  = "https://github.com/OpenLocalizationTest/oltest/blob/002d5baaaa06fb18bf2a97ca5f648b13fd56d8a5/e2e/"
$zhcnUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/15ecd35071fdc19de379c61081726885a176079d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$dedeUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/be25bba22b2d8fab59d6d9783921ba3b0a3f9129/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

# =======================================================================
# Sheet "Overview" — one row per source file
# =======================================================================
$ov = $wb.Worksheets.Item("Overview")

# Dropping the stale hyperlink now; every hyperlink on the sheet gets
# re-added below in the order it must appear (A2, A3, A4).
$ov.Range("A2").Hyperlinks.Delete()

# Row 2 now describes the first .png dependency (previously the .md file)
$ov.Range("A2").Value2 = $png1Name
$ov.Range("B2").Value2 = $readyStatus
$ov.Range("C2").Value2 = $readyStatus
$ov.Range("D2").Value2 = $overviewDate
Add-StyledHyperlink $ov "A2" ($srcUrlBase + $png1Name) $png1Name

# Row 3 — second .png dependency
$ov.Range("B3").Value2 = $readyStatus
$ov.Range("C3").Value2 = $readyStatus
$ov.Range("D3").Value2 = $overviewDate
Add-StyledHyperlink $ov "A3" ($srcUrlBase + $png2Name) $png2Name

# Row 4 — the .md file itself
$ov.Range("B4").Value2 = $readyStatus
$ov.Range("C4").Value2 = $readyStatus
$ov.Range("D4").Value2 = $overviewDate
Add-StyledHyperlink $ov "A4" ($srcUrlBase + $mdName) $mdName

# =======================================================================
# Per-locale detail sheets (zh-cn, de-de)
# =======================================================================
$locales = @(
    @{ Sheet = "zh-cn"; UrlBase = $zhcnUrlBase; TargetMd = $mdTargetZhCn; HandoffDate = "2016-03-20 16:55:19" },
    @{ Sheet = "de-de"; UrlBase = $dedeUrlBase; TargetMd = $mdTargetDeDe; HandoffDate = "2016-03-20 16:55:23" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # Drop the stale hyperlinks; all 9 needed ones (A/B/D x rows 2-4) are
    # re-added below in that exact order.
    $ws.Range("A2").Hyperlinks.Delete()

    # ---- Row 2: first .png dependency ---------------------------------
    $ws.Range("A2").Value2 = $png1Name
    $ws.Range("B2").Value2 = ".png"
    $ws.Range("C2").Value2 = $readyStatus
    $ws.Range("D2").Value2 = $png1TargetName
    $ws.Range("E2").Value2 = $loc.HandoffDate
    $ws.Range("H2").Value2 = "0001-01-01 00:00:00"
    $ws.Range("I2").Value2 = "IsDependency"
    $ws.Range("J2").Value2 = "e2e\" + $mdName

    Add-StyledHyperlink $ws "A2" ($srcUrlBase + $png1Name) $png1Name
    Add-StyledHyperlink $ws "B2" ($srcUrlBase + $png1Name) ".png"
    Add-StyledHyperlink $ws "D2" ($loc.UrlBase + $png1TargetName) $png1TargetName

    # ---- Row 3: second .png dependency ---------------------------------
    $ws.Range("A3").Value2 = $png2Name
    $ws.Range("B3").Value2 = ".png"
    $ws.Range("C3").Value2 = $readyStatus
    $ws.Range("D3").Value2 = $png2TargetName
    $ws.Range("E3").Value2 = $loc.HandoffDate
    $ws.Range("H3").Value2 = "0001-01-01 00:00:00"
    $ws.Range("I3").Value2 = "IsDependency"
    $ws.Range("J3").Value2 = "e2e\" + $mdName

    Add-StyledHyperlink $ws "A3" ($srcUrlBase + $png2Name) $png2Name
    Add-StyledHyperlink $ws "B3" ($srcUrlBase + $png2Name) ".png"
    Add-StyledHyperlink $ws "D3" ($loc.UrlBase + $png2TargetName) $png2TargetName

    # ---- Row 4: the .md file itself ------------------------------------
    $ws.Range("A4").Value2 = $mdName
    $ws.Range("B4").Value2 = ".md"
    $ws.Range("C4").Value2 = $readyStatus
    $ws.Range("D4").Value2 = $loc.TargetMd
    $ws.Range("E4").Value2 = $loc.HandoffDate
    $ws.Range("H4").Value2 = "0001-01-01 00:00:00"
    $ws.Range("I4").Value2 = "Include"

    Add-StyledHyperlink $ws "A4" ($srcUrlBase + $mdName) $mdName
    Add-StyledHyperlink $ws "B4" ($srcUrlBase + $mdName) ".md"
    Add-StyledHyperlink $ws "D4" ($loc.UrlBase + $loc.TargetMd) $loc.TargetMd

    # Column E keeps the workbook's custom date-time number format.
    $ws.Range("E2:E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

Write-Output "Report regenerated for handoff."
